$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Deals_data")

# Update existing row 2: providerName/dealName change from Peeq -> FMTLocal deal
$ws1.Range("B2").Value = "FMTLocal"
$ws1.Range("C2").Value = "FmtLocalEssentialsDeal_Safebase1_Bundle_DealDescription"

# Add new row 3: Driven deal (VAS)
$ws1.Range("A3").Value = "VAS"
$ws1.Range("B3").Value = "Driven"
$ws1.Range("C3").Value = "Driven @ R139"
$ws1.Range("D3").Value = "'0678678769"
$ws1.Range("E3").Value = "Durban"

# Add new row 4: On Air deal (Telco)
$ws1.Range("A4").Value = "Telco"
$ws1.Range("B4").Value = "On Air"
$ws1.Range("C4").Value = "On Air Testing Deal"
$ws1.Range("D4").Value = "'0678678770"
$ws1.Range("E4").Value = "Durban"

# Update the selection on the reference "Sheet1" tab to match the new data extent
$ws3 = $wb.Worksheets.Item("Sheet1")
$ws3.Range("A2:E4").Select()

# Re-activate Deals_data and select the full updated range, restoring it as the active tab
$ws1.Activate()
$ws1.Range("A2:E4").Select()
